$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns at F:G (old F..M shift to H..O).
$ws.Range("F1:G1").EntireColumn.Insert()

# Insert two more columns further out (still inside the originally-uniform-width
# zone) so the column-width block correctly extends to cover the 4 new columns
# without needing to hand-set any ColumnWidth values (which would quantize the
# stored width away from the exact original 16.7109375).
$ws.Range("U1:V1").EntireColumn.Insert()

# --- Row 1 (headers) ---
$ws.Range("F1").Value = "Most aligned target"
$ws.Range("G1").Value = "Least aligned target"
$ws.Range("M1").Value = "Alignment score (Delegator|ADM(most))"
$ws.Range("N1").Value = "Alignment score (Delegator|ADM(least))"
$ws.Range("O1").Value = "Match_MostAligned"
$ws.Range("P1").Value = "Mach_LeastAligned"
$ws.Range("Q1").Value = "Match_GrpMembers"

# --- Row 2 (source) ---
# F2/G2 have no content in the new layout - remove the placeholder cells that
# the column insert created (they picked up neighboring style/format).
$ws.Range("F2:G2").Clear()
$ws.Range("M2").Value = "TA1 server"
$ws.Range("N2").Value = "TA1 server"
$ws.Range("O2").Value = "Calculated from probe responses"
$ws.Range("P2").Value = "Calculated from probe responses"
$ws.Range("Q2").Value = "Calculated from probe responses"

# --- Row 3 (definition) ---
$ws.Range("F3").Value = "Target with the highest alignment score for the delegator on the text scenario"
$ws.Range("G3").Value = "Target with the lowest alignment score for the delegator on the text scenario"
$ws.Range("H3").Value = "Calculated alignment score between the delegator and a target"
$ws.Range("I3").Value = "Calculated alignment score between the delegator and a target"
$ws.Range("M3").Value = "Calculated alignment score between the delegator and the aligned ADM run on the same scenario at the most aligned target"
$ws.Range("N3").Value = "Calculated alignment score between the delegator and the aligned ADM run on the same scenario at the least aligned target"
$ws.Range("O3").Value = "% of exact matches on probe responses between delegator and ADM run on same scenario at most aligned target"
$ws.Range("P3").Value = "% of exact matches on probe responses between delegator and ADM run on same scenario at least aligned target"
$ws.Range("Q3").Value = "% of exact matches on probe responses between delegator and ADM run on same scenario at group target"

# --- Row 4 (levels) ---
# F4/G4 have no content in the new layout - remove the placeholder cells.
$ws.Range("F4:G4").Clear()
# M4/N4 also have no content in the new layout.
$ws.Range("M4:N4").Clear()
$ws.Range("O4").Value = "0-100"
$ws.Range("P4").Value = "0-100"
$ws.Range("Q4").Value = "0-100"

# Row 3 grew taller to fit the newly added, longer definition text.
$ws.Range("A3").RowHeight = 137.25

# Update the active selection to match the saved view state.
$ws.Range("N4").Select()
